# Insert a new weekly price record for Zanahoria (Terminal Hortofrutícola Agro
# Chillán) as row 381, shifting the existing rows 381-407 down to 382-408.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 381 (pushes rows 381..407 to 382..408).
$ws.Rows.Item(381).Insert()

# Populate the new row 381 with the new record.
$ws.Cells.Item(381, 1).Value  = 7
$ws.Cells.Item(381, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(381, 3).Value  = "Ñuble"
$ws.Cells.Item(381, 4).Value  = Get-Date -Year 2023 -Month 4 -Day 5 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(381, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(381, 5).Value  = 16
$ws.Cells.Item(381, 6).Value  = 100114013
$ws.Cells.Item(381, 7).Value  = "Zanahoria"
$ws.Cells.Item(381, 8).Value  = "Sin especificar"
$ws.Cells.Item(381, 9).Value  = "Primera"
$ws.Cells.Item(381, 10).Value = 80
$ws.Cells.Item(381, 11).Value = 6000
$ws.Cells.Item(381, 12).Value = 7000
$ws.Cells.Item(381, 13).Value = 6500
$ws.Cells.Item(381, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(381, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(381, 16).Value = 325
$ws.Cells.Item(381, 17).Value = 20
$ws.Cells.Item(381, 18).Value = "Hortaliza"
